$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo in header: MidleName -> MiddleName
$ws.Range("C1").Value = "MiddleName"

# Update Classroom (column U) values for rows 2-100
$ws.Range("U2").Value = 0
$ws.Range("U3").Value = 3
$ws.Range("U4").Value = 2
$ws.Range("U5").Value = 1
$ws.Range("U6").Value = 0
$ws.Range("U7").Value = 3
$ws.Range("U8").Value = 1
$ws.Range("U9").Value = 0
$ws.Range("U10").Value = 2
$ws.Range("U11").Value = 3
$ws.Range("U12").Value = 2
$ws.Range("U13").Value = 1
$ws.Range("U14").Value = 0
$ws.Range("U15").Value = 1
$ws.Range("U16").Value = 3
$ws.Range("U17").Value = 1
$ws.Range("U18").Value = 2
$ws.Range("U19").Value = 1
$ws.Range("U20").Value = 3
$ws.Range("U21").Value = 2
$ws.Range("U22").Value = 0
$ws.Range("U23").Value = 1
$ws.Range("U24").Value = 2
$ws.Range("U25").Value = 1
$ws.Range("U26").Value = 0
$ws.Range("U27").Value = 3
$ws.Range("U28").Value = 1
$ws.Range("U29").Value = 0
$ws.Range("U30").Value = 3
$ws.Range("U31").Value = 2
$ws.Range("U32").Value = 0
$ws.Range("U33").Value = 3
$ws.Range("U34").Value = 3
$ws.Range("U35").Value = 0
$ws.Range("U36").Value = 1
$ws.Range("U37").Value = 2
$ws.Range("U38").Value = 3
$ws.Range("U39").Value = 3
$ws.Range("U40").Value = 1
$ws.Range("U41").Value = 3
$ws.Range("U42").Value = 0
$ws.Range("U43").Value = 2
$ws.Range("U44").Value = 1
$ws.Range("U45").Value = 2
$ws.Range("U46").Value = 0
$ws.Range("U47").Value = 2
$ws.Range("U48").Value = 0
$ws.Range("U49").Value = 3
$ws.Range("U50").Value = 0
$ws.Range("U51").Value = 2
$ws.Range("U52").Value = 1
$ws.Range("U53").Value = 0
$ws.Range("U54").Value = 2
$ws.Range("U55").Value = 3
$ws.Range("U56").Value = 1
$ws.Range("U57").Value = 2
$ws.Range("U58").Value = 2
$ws.Range("U59").Value = 3
$ws.Range("U60").Value = 0
$ws.Range("U61").Value = 1
$ws.Range("U62").Value = 0
$ws.Range("U63").Value = 2
$ws.Range("U64").Value = 1
$ws.Range("U65").Value = 2
$ws.Range("U66").Value = 3
$ws.Range("U67").Value = 3
$ws.Range("U68").Value = 0
$ws.Range("U69").Value = 3
$ws.Range("U70").Value = 0
$ws.Range("U71").Value = 2
$ws.Range("U72").Value = 3
$ws.Range("U73").Value = 2
$ws.Range("U74").Value = 0
$ws.Range("U75").Value = 2
$ws.Range("U76").Value = 1
$ws.Range("U77").Value = 1
$ws.Range("U78").Value = 0
$ws.Range("U79").Value = 1
$ws.Range("U80").Value = 3
$ws.Range("U81").Value = 2
$ws.Range("U82").Value = 3
$ws.Range("U83").Value = 1
$ws.Range("U84").Value = 1
$ws.Range("U85").Value = 0
$ws.Range("U86").Value = 2
$ws.Range("U87").Value = 1
$ws.Range("U88").Value = 0
$ws.Range("U89").Value = 1
$ws.Range("U90").Value = 3
$ws.Range("U91").Value = 0
$ws.Range("U92").Value = 2
$ws.Range("U93").Value = 3
$ws.Range("U94").Value = 1
$ws.Range("U95").Value = 0
$ws.Range("U96").Value = 2
$ws.Range("U97").Value = 1
$ws.Range("U98").Value = 3
$ws.Range("U99").Value = 0
$ws.Range("U100").Value = 2
